$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad" / "Changed") holds a date serial that is bumped by
# one day on every automatic refresh of this report (rows 2-433 all share
# the same value). Increment each of them from 46075 to 46076.
for ($r = 2; $r -le 433; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
